{"js": "// Replace each two-digit multiplication expression in the practice\n// table with its updated version. Matches are unique, exact strings\n// (e.g. \"68\u00d763=\") so a simple body.search + insertText(replace) per\n// pair is unambiguous.\nconst replacements = [[\"68\u00d763=\", \"99\u00d773=\"], [\"65\u00d770=\", \"16\u00d796=\"], [\"58\u00d741=\", \"60\u00d781=\"], [\"96\u00d713=\", \"77\u00d783=\"], [\"11\u00d714=\", \"23\u00d786=\"], [\"85\u00d746=\", \"56\u00d799=\"], [\"38\u00d750=\", \"48\u00d741=\"], [\"68\u00d759=\", \"70\u00d762=\"], [\"22\u00d771=\", \"22\u00d756=\"], [\"10\u00d775=\", \"39\u00d783=\"], [\"82\u00d726=\", \"17\u00d772=\"], [\"33\u00d740=\", \"37\u00d799=\"], [\"94\u00d734=\", \"47\u00d773=\"], [\"68\u00d761=\", \"79\u00d715=\"], [\"23\u00d736=\", \"57\u00d726=\"], [\"78\u00d772=\", \"74\u00d762=\"], [\"63\u00d774=\", \"82\u00d776=\"], [\"19\u00d794=\", \"21\u00d719=\"], [\"93\u00d790=\", \"63\u00d799=\"], [\"29\u00d796=\", \"13\u00d7100=\"], [\"18\u00d742=\", \"74\u00d736=\"], [\"67\u00d788=\", \"28\u00d730=\"], [\"41\u00d760=\", \"44\u00d795=\"], [\"90\u00d752=\", \"85\u00d762=\"], [\"98\u00d752=\", \"100\u00d745=\"], [\"53\u00d772=\", \"87\u00d744=\"], [\"94\u00d775=\", \"17\u00d767=\"], [\"75\u00d752=\", \"77\u00d728=\"], [\"80\u00d799=\", \"94\u00d721=\"], [\"38\u00d732=\", \"94\u00d715=\"], [\"28\u00d728=\", \"63\u00d784=\"], [\"89\u00d777=\", \"68\u00d789=\"], [\"77\u00d767=\", \"10\u00d710=\"], [\"46\u00d727=\", \"33\u00d738=\"], [\"64\u00d738=\", \"31\u00d713=\"], [\"30\u00d767=\", \"33\u00d796=\"], [\"70\u00d790=\", \"29\u00d779=\"], [\"52\u00d758=\", \"100\u00d757=\"], [\"80\u00d772=\", \"62\u00d734=\"], [\"75\u00d744=\", \"28\u00d729=\"], [\"23\u00d713=\", \"64\u00d753=\"], [\"46\u00d780=\", \"20\u00d730=\"], [\"63\u00d710=\", \"85\u00d719=\"], [\"70\u00d787=\", \"69\u00d710=\"], [\"71\u00d743=\", \"82\u00d799=\"], [\"54\u00d766=\", \"56\u00d774=\"], [\"51\u00d788=\", \"62\u00d739=\"], [\"89\u00d747=\", \"15\u00d765=\"], [\"68\u00d778=\", \"75\u00d792=\"], [\"36\u00d726=\", \"28\u00d790=\"], [\"56\u00d734=\", \"42\u00d760=\"], [\"32\u00d744=\", \"57\u00d7100=\"], [\"77\u00d735=\", \"83\u00d710=\"], [\"24\u00d755=\", \"49\u00d794=\"], [\"54\u00d732=\", \"42\u00d721=\"], [\"92\u00d739=\", \"46\u00d776=\"], [\"46\u00d785=\", \"87\u00d784=\"], [\"19\u00d734=\", \"95\u00d737=\"], [\"39\u00d728=\", \"72\u00d799=\"], [\"44\u00d738=\", \"46\u00d729=\"], [\"85\u00d724=\", \"32\u00d719=\"], [\"91\u00d798=\", \"92\u00d715=\"], [\"62\u00d792=\", \"25\u00d713=\"], [\"46\u00d758=\", \"51\u00d733=\"], [\"19\u00d791=\", \"16\u00d799=\"], [\"60\u00d710=\", \"98\u00d779=\"], [\"52\u00d785=\", \"87\u00d797=\"], [\"45\u00d753=\", \"24\u00d721=\"], [\"63\u00d720=\", \"76\u00d718=\"], [\"95\u00d754=\", \"93\u00d777=\"], [\"77\u00d771=\", \"47\u00d737=\"], [\"68\u00d780=\", \"26\u00d798=\"], [\"41\u00d758=\", \"70\u00d738=\"], [\"16\u00d759=\", \"59\u00d710=\"], [\"45\u00d741=\", \"46\u00d716=\"], [\"83\u00d774=\", \"83\u00d759=\"], [\"50\u00d769=\", \"31\u00d737=\"], [\"82\u00d725=\", \"100\u00d793=\"], [\"16\u00d757=\", \"21\u00d726=\"], [\"53\u00d722=\", \"81\u00d755=\"], [\"80\u00d760=\", \"21\u00d789=\"], [\"62\u00d787=\", \"53\u00d737=\"], [\"38\u00d762=\", \"20\u00d785=\"], [\"65\u00d787=\", \"34\u00d763=\"], [\"73\u00d798=\", \"30\u00d789=\"], [\"40\u00d723=\", \"74\u00d720=\"], [\"83\u00d725=\", \"72\u00d712=\"], [\"32\u00d764=\", \"99\u00d752=\"], [\"10\u00d756=\", \"35\u00d740=\"], [\"55\u00d722=\", \"14\u00d714=\"], [\"22\u00d7100=\", \"46\u00d726=\"], [\"19\u00d774=\", \"36\u00d779=\"], [\"43\u00d730=\", \"33\u00d765=\"], [\"10\u00d729=\", \"65\u00d755=\"], [\"42\u00d731=\", \"33\u00d732=\"], [\"25\u00d748=\", \"82\u00d734=\"], [\"49\u00d788=\", \"97\u00d762=\"], [\"88\u00d743=\", \"14\u00d734=\"], [\"35\u00d717=\", \"68\u00d727=\"], [\"70\u00d797=\", \"50\u00d793=\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each two-digit multiplication expression in the practice\n# table to its new value. Each old string is unique in the document,\n# so Find/Execute with MatchCase + MatchWholeWord=False safely\n# targets exactly one run per pair (Wrap=wdFindContinue=1, Replace=wdReplaceAll=2).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"68\u00d763=\", \"99\u00d773=\"),\n  @(\"65\u00d770=\", \"16\u00d796=\"),\n  @(\"58\u00d741=\", \"60\u00d781=\"),\n  @(\"96\u00d713=\", \"77\u00d783=\"),\n  @(\"11\u00d714=\", \"23\u00d786=\"),\n  @(\"85\u00d746=\", \"56\u00d799=\"),\n  @(\"38\u00d750=\", \"48\u00d741=\"),\n  @(\"68\u00d759=\", \"70\u00d762=\"),\n  @(\"22\u00d771=\", \"22\u00d756=\"),\n  @(\"10\u00d775=\", \"39\u00d783=\"),\n  @(\"82\u00d726=\", \"17\u00d772=\"),\n  @(\"33\u00d740=\", \"37\u00d799=\"),\n  @(\"94\u00d734=\", \"47\u00d773=\"),\n  @(\"68\u00d761=\", \"79\u00d715=\"),\n  @(\"23\u00d736=\", \"57\u00d726=\"),\n  @(\"78\u00d772=\", \"74\u00d762=\"),\n  @(\"63\u00d774=\", \"82\u00d776=\"),\n  @(\"19\u00d794=\", \"21\u00d719=\"),\n  @(\"93\u00d790=\", \"63\u00d799=\"),\n  @(\"29\u00d796=\", \"13\u00d7100=\"),\n  @(\"18\u00d742=\", \"74\u00d736=\"),\n  @(\"67\u00d788=\", \"28\u00d730=\"),\n  @(\"41\u00d760=\", \"44\u00d795=\"),\n  @(\"90\u00d752=\", \"85\u00d762=\"),\n  @(\"98\u00d752=\", \"100\u00d745=\"),\n  @(\"53\u00d772=\", \"87\u00d744=\"),\n  @(\"94\u00d775=\", \"17\u00d767=\"),\n  @(\"75\u00d752=\", \"77\u00d728=\"),\n  @(\"80\u00d799=\", \"94\u00d721=\"),\n  @(\"38\u00d732=\", \"94\u00d715=\"),\n  @(\"28\u00d728=\", \"63\u00d784=\"),\n  @(\"89\u00d777=\", \"68\u00d789=\"),\n  @(\"77\u00d767=\", \"10\u00d710=\"),\n  @(\"46\u00d727=\", \"33\u00d738=\"),\n  @(\"64\u00d738=\", \"31\u00d713=\"),\n  @(\"30\u00d767=\", \"33\u00d796=\"),\n  @(\"70\u00d790=\", \"29\u00d779=\"),\n  @(\"52\u00d758=\", \"100\u00d757=\"),\n  @(\"80\u00d772=\", \"62\u00d734=\"),\n  @(\"75\u00d744=\", \"28\u00d729=\"),\n  @(\"23\u00d713=\", \"64\u00d753=\"),\n  @(\"46\u00d780=\", \"20\u00d730=\"),\n  @(\"63\u00d710=\", \"85\u00d719=\"),\n  @(\"70\u00d787=\", \"69\u00d710=\"),\n  @(\"71\u00d743=\", \"82\u00d799=\"),\n  @(\"54\u00d766=\", \"56\u00d774=\"),\n  @(\"51\u00d788=\", \"62\u00d739=\"),\n  @(\"89\u00d747=\", \"15\u00d765=\"),\n  @(\"68\u00d778=\", \"75\u00d792=\"),\n  @(\"36\u00d726=\", \"28\u00d790=\"),\n  @(\"56\u00d734=\", \"42\u00d760=\"),\n  @(\"32\u00d744=\", \"57\u00d7100=\"),\n  @(\"77\u00d735=\", \"83\u00d710=\"),\n  @(\"24\u00d755=\", \"49\u00d794=\"),\n  @(\"54\u00d732=\", \"42\u00d721=\"),\n  @(\"92\u00d739=\", \"46\u00d776=\"),\n  @(\"46\u00d785=\", \"87\u00d784=\"),\n  @(\"19\u00d734=\", \"95\u00d737=\"),\n  @(\"39\u00d728=\", \"72\u00d799=\"),\n  @(\"44\u00d738=\", \"46\u00d729=\"),\n  @(\"85\u00d724=\", \"32\u00d719=\"),\n  @(\"91\u00d798=\", \"92\u00d715=\"),\n  @(\"62\u00d792=\", \"25\u00d713=\"),\n  @(\"46\u00d758=\", \"51\u00d733=\"),\n  @(\"19\u00d791=\", \"16\u00d799=\"),\n  @(\"60\u00d710=\", \"98\u00d779=\"),\n  @(\"52\u00d785=\", \"87\u00d797=\"),\n  @(\"45\u00d753=\", \"24\u00d721=\"),\n  @(\"63\u00d720=\", \"76\u00d718=\"),\n  @(\"95\u00d754=\", \"93\u00d777=\"),\n  @(\"77\u00d771=\", \"47\u00d737=\"),\n  @(\"68\u00d780=\", \"26\u00d798=\"),\n  @(\"41\u00d758=\", \"70\u00d738=\"),\n  @(\"16\u00d759=\", \"59\u00d710=\"),\n  @(\"45\u00d741=\", \"46\u00d716=\"),\n  @(\"83\u00d774=\", \"83\u00d759=\"),\n  @(\"50\u00d769=\", \"31\u00d737=\"),\n  @(\"82\u00d725=\", \"100\u00d793=\"),\n  @(\"16\u00d757=\", \"21\u00d726=\"),\n  @(\"53\u00d722=\", \"81\u00d755=\"),\n  @(\"80\u00d760=\", \"21\u00d789=\"),\n  @(\"62\u00d787=\", \"53\u00d737=\"),\n  @(\"38\u00d762=\", \"20\u00d785=\"),\n  @(\"65\u00d787=\", \"34\u00d763=\"),\n  @(\"73\u00d798=\", \"30\u00d789=\"),\n  @(\"40\u00d723=\", \"74\u00d720=\"),\n  @(\"83\u00d725=\", \"72\u00d712=\"),\n  @(\"32\u00d764=\", \"99\u00d752=\"),\n  @(\"10\u00d756=\", \"35\u00d740=\"),\n  @(\"55\u00d722=\", \"14\u00d714=\"),\n  @(\"22\u00d7100=\", \"46\u00d726=\"),\n  @(\"19\u00d774=\", \"36\u00d779=\"),\n  @(\"43\u00d730=\", \"33\u00d765=\"),\n  @(\"10\u00d729=\", \"65\u00d755=\"),\n  @(\"42\u00d731=\", \"33\u00d732=\"),\n  @(\"25\u00d748=\", \"82\u00d734=\"),\n  @(\"49\u00d788=\", \"97\u00d762=\"),\n  @(\"88\u00d743=\", \"14\u00d734=\"),\n  @(\"35\u00d717=\", \"68\u00d727=\"),\n  @(\"70\u00d797=\", \"50\u00d793=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $found = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    throw \"No match found for: $old\"\n  }\n}\n"}
